$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 122 (old row 122 data shifts down to row 124)
$ws.Rows.Item(122).Insert()
$ws.Rows.Item(122).Insert()

# --- New row 122 ---
$ws.Cells.Item(122, 1).Value = 7
$ws.Cells.Item(122, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(122, 3).Value = "Ñuble"
$ws.Cells.Item(122, 4).Value = 44448
$ws.Cells.Item(122, 4).Style = $ws.Cells.Item(121, 4).Style
$ws.Cells.Item(122, 4).NumberFormat = $ws.Cells.Item(121, 4).NumberFormat
$ws.Cells.Item(122, 5).Value = 16
$ws.Cells.Item(122, 6).Value = "Fruta"
$ws.Cells.Item(122, 7).Value = 100108
$ws.Cells.Item(122, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(122, 9).Value = 100108005
$ws.Cells.Item(122, 10).Value = "Piña"
$ws.Cells.Item(122, 11).Value = "Caramelo"
$ws.Cells.Item(122, 12).Value = "Primera"
$ws.Cells.Item(122, 13).Value = 60
$ws.Cells.Item(122, 14).Value = 19000
$ws.Cells.Item(122, 15).Value = 20000
$ws.Cells.Item(122, 16).Value = 19500
$ws.Cells.Item(122, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(122, 18).Value = "Ecuador"
$ws.Cells.Item(122, 19).Value = 1625
$ws.Cells.Item(122, 20).Value = 12

# --- New row 123 ---
$ws.Cells.Item(123, 1).Value = 7
$ws.Cells.Item(123, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(123, 3).Value = "Ñuble"
$ws.Cells.Item(123, 4).Value = 44448
$ws.Cells.Item(123, 4).Style = $ws.Cells.Item(121, 4).Style
$ws.Cells.Item(123, 4).NumberFormat = $ws.Cells.Item(121, 4).NumberFormat
$ws.Cells.Item(123, 5).Value = 16
$ws.Cells.Item(123, 6).Value = "Fruta"
$ws.Cells.Item(123, 7).Value = 100108
$ws.Cells.Item(123, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(123, 9).Value = 100108005
$ws.Cells.Item(123, 10).Value = "Piña"
$ws.Cells.Item(123, 11).Value = "Caramelo"
$ws.Cells.Item(123, 12).Value = "Segunda"
$ws.Cells.Item(123, 13).Value = 60
$ws.Cells.Item(123, 14).Value = 19000
$ws.Cells.Item(123, 15).Value = 20000
$ws.Cells.Item(123, 16).Value = 19500
$ws.Cells.Item(123, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(123, 18).Value = "Ecuador"
$ws.Cells.Item(123, 19).Value = 1393
$ws.Cells.Item(123, 20).Value = 14

$wb.Save()
